# Update "想去人数" (column F) values on the "展览" and "全部类型" sheets.
# Both sheets hold the same event list, so the same row/value updates apply
# to each.

$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 1096
    5  = 4668
    7  = 395
    8  = 1401
    10 = 56
    11 = 1151
    13 = 637
    15 = 40
    16 = 18
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
